# Replace the Python dict-literal text in A1 with a pretty-printed JSON
# equivalent (the smart quotes around "Hello" become the literal \u2018/\u2019
# escape sequences that json.dumps(..., ensure_ascii=True) produces), move it
# from A2 (shared-string duplicate) up into A1, drop the old numeric A1 value,
# and strip the bold/bordered/centered formatting that used to mark A1 as a
# "header" cell so the sheet ends up as a single, unformatted A1 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = 'questions = [
    {
        "title": "Considering the code below, how many times will the \u2018Hello\u2019 message be displayed on the console? const App = (props) =&gt {\n const [counter, setCounter] = useState(0)\n useEffect(\n   () =&gt {\n     console.log(''Hello'')\n     setCounter(1)\n   },\n   [props.visible]\n )\n return &ltdiv&gt{counter}&lt/div&gt\n}",
        "ques_type": 2,
        "options": [
            "0",
            "1",
            "2",
            "3"
        ],
        "score": "1"
    },
    {
        "title": "Which statement describes the code below? const fetchData = () =&gt new Promise((r) =&gt setTimeout(() =&gt r(Date.now()), 100))\n \nconst MyComponent = () =&gt {\n const [result, setResult] = React.useState()\n const data = fetchData().then((value) =&gt setResult(value))\n return (\n   &ltdiv&gt\n     {result === data.toString() ? (\n       &ltdiv&gthello&lt/div&gt\n     ) : (\n       &ltdiv&gtgood bye&lt/div&gt\n     )}\n   &lt/div&gt\n )\n}",
        "ques_type": 2,
        "options": [
            "A ''good bye'' message will be displayed.",
            "A ''hello'' message will be displayed.",
            "The code results in a memory leak.",
            "setResult is never called."
        ],
        "score": "A ''good bye'' message will be displayed."
    },
    {
        "title": "Which wrapper will hide its child component for four seconds?",
        "ques_type": 2,
        "options": [
            "const HiderWrapper = (props) =&gt {\n return (\n   setTimeout(() =&gt {\n     return props.children\n   }),\n   4000\n )\n}\n",
            "const HiderWrapper = (props) =&gt {\n const [visible, setVisible] = useState(false)\n useEffect(() =&gt {\n   setInterval(() =&gt {\n     setVisible(true)\n   }, 4000)\n }, [])\n if (visible) return props.body\n else return null\n}\n",
            "const HiderWrapper = (props) =&gt {\n const [visible, setVisible] = useState(false)\n useEffect(() =&gt {\n   setTimeout(() =&gt {\n     setVisible(true)\n   }, 4000)\n }, [])\n if (visible) return props.children\n else return null\n}\n",
            "const HiderWrapper = () =&gt {\n wait(4)\n return null\n}\n"
        ],
        "score": "const HiderWrapper = (props) =&gt {\n const [visible, setVisible] = useState(false)\n useEffect(() =&gt {\n   setTimeout(() =&gt {\n     setVisible(true)\n   }, 4000)\n }, [])\n if (visible) return props.children\n else return null\n}"
    },
    {
        "title": "Considering the code below, when will the MyChild component be unmounted? const MyParent = ({ value }) =&gt {\n return &ltdiv&gt{value !== 3 &amp&amp &ltMyChild /&gt}&lt/div&gt\n}",
        "ques_type": 2,
        "options": [
            "When the value property is equal to 3.",
            "When the value property is different from 3.",
            "Never.",
            "After each render of the MyParent component."
        ],
        "score": "When the value property is equal to 3."
    }
]'

# A2 held the shared-string text previously; clear it first so the old
# (pre-reformat) shared string has no remaining references once A1 is
# rewritten below.
$ws.Range("A2").ClearContents()
$ws.Range("A2").ClearFormats()

# A1 held a bold/bordered/centered numeric 0 - strip that formatting so the
# cell reverts to the default (unstyled) look.
$ws.Range("A1").ClearFormats()

# Write the reformatted text into A1.
$ws.Range("A1").Value = $newText

# The new text contains embedded newlines; re-run autofit so the row height
# stays at the sheet's default instead of being stretched to fit every line.
$ws.Range("A1").EntireRow.AutoFit()
